$d = $word.ActiveDocument

# The SUS ("System Usability Scale") scoring table ends with a
# "Valutazione totale" (total score) row whose last cell holds the
# overall score. Locate that table, then update the score cell from
# "19" to "27,5".
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text -like "*Valutazione totale*") {
        $table = $candidate
    }
}

$lastRow = $table.Rows.Count
$lastCol = $table.Columns.Count
$cell = $table.Cell($lastRow, $lastCol)

$cellRange = $cell.Range
if ($cellRange.Text.TrimEnd([char]13, [char]7) -eq "19") {
    $cellRange.Text = "27,5"
}
